$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# New test cases for Use Case "7, Search" (sheet1) + matching test-run log
# rows (sheet2). Values are written in the same interleaved order the
# original author used so shared-string indices line up.
# ---------------------------------------------------------------------------

# --- sheet1 row 17 / sheet2 row 17 -----------------------------------------
$ws1.Cells.Item(17,2).Value = "7, Search"
$ws1.Cells.Item(17,1).Value = 16
$ws1.Cells.Item(17,3).Value = "Main Flow"
$ws1.Cells.Item(17,4).Value = "Benutzerinput pdf und mit OK bestätigen"
$ws1.Cells.Item(17,5).Value = "Test database(Dokumente .pdf, cdew.pdf und test2009.pdf vorhanden)"
$ws1.Cells.Item(17,6).Value = "Meldung mit Liste der drei Dokumente .pdf, cdew.pdf und test2009.pdf"
$ws1.Range("E2:F2").Copy()
$ws1.Range("E17:F17").PasteSpecial(-4122)
$ws1.Rows.Item(17).RowHeight = 47.25

$ws2.Cells.Item(17,1).Value = 45132
$ws2.Cells.Item(17,2).Value = 0.68819444444444444
$ws2.Cells.Item(17,3).Value = "main"
$ws2.Range("A2:C2").Copy()
$ws2.Range("A17:C17").PasteSpecial(-4122)
$ws2.Cells.Item(17,4).Value = 16
$ws2.Cells.Item(17,5).Value = "Kierstein"
$ws2.Cells.Item(17,6).Value = "Success"

# --- sheet1 row 18 / sheet2 row 18 -----------------------------------------
$ws1.Cells.Item(18,1).Value = 17
$ws1.Cells.Item(18,2).Value = "7, Search"
$ws1.Cells.Item(18,3).Value = "Alternative Flow 1"
$ws1.Cells.Item(18,4).Value = "Bei Abfrage auf Input abbrechen wählen"
$ws1.Cells.Item(18,5).Value = "Test database()"
$ws1.Cells.Item(18,6).Value = "Benutzer landet wieder im Menü"
$ws1.Range("E2:F2").Copy()
$ws1.Range("E18:F18").PasteSpecial(-4122)

$ws2.Cells.Item(18,1).Value = 45132
$ws2.Cells.Item(18,2).Value = 0.68958333333333333
$ws2.Cells.Item(18,3).Value = "main"
$ws2.Range("A2:C2").Copy()
$ws2.Range("A18:C18").PasteSpecial(-4122)
$ws2.Cells.Item(18,4).Value = 17
$ws2.Cells.Item(18,5).Value = "Kierstein"
$ws2.Cells.Item(18,6).Value = "Success"

# --- sheet1 row 19 / sheet2 row 19 -----------------------------------------
$ws1.Cells.Item(19,1).Value = 18
$ws1.Cells.Item(19,2).Value = "7, Search"
$ws1.Cells.Item(19,3).Value = "Alternative Flow 2"
$ws1.Cells.Item(19,4).Value = "Benutzerinput PDF und mit OK bestätigen"
$ws1.Cells.Item(19,5).Value = "Test database(Ohne Dokumente mit PDF im Namen)"
$ws1.Cells.Item(19,6).Value = "Da keine Dokumente gefunden wurden soll der Benutzer erneut Input geben"
$ws1.Range("E2:F2").Copy()
$ws1.Range("E19:F19").PasteSpecial(-4122)
$ws1.Cells.Item(19,7).Value = "Test für Extension 3a nach Cockburn Template"
$ws1.Rows.Item(19).RowHeight = 47.25

$ws2.Cells.Item(19,1).Value = 45132
$ws2.Cells.Item(19,2).Value = 0.69166666666666676
$ws2.Cells.Item(19,3).Value = "main"
$ws2.Range("A2:C2").Copy()
$ws2.Range("A19:C19").PasteSpecial(-4122)
$ws2.Cells.Item(19,4).Value = 18
$ws2.Cells.Item(19,5).Value = "Kierstein"
$ws2.Cells.Item(19,6).Value = "Failed"
$ws2.Cells.Item(19,7).Value = "Der Benutzer erhält eine leere Meldung und landet wieder im Menü anstatt direkt neuen Input geben zu können"

# --- sheet1 row 20 / sheet2 row 20 -----------------------------------------
$ws1.Cells.Item(20,1).Value = 19
$ws1.Cells.Item(20,2).Value = "7, Search"
$ws1.Cells.Item(20,3).Value = "Alternative Flow 3"
$ws1.Cells.Item(20,7).Value = "Test für Extension 3b nach Cockburn Template"

$ws2.Cells.Item(20,1).Value = 45132
$ws2.Cells.Item(20,2).Value = 0.69305555555555554
$ws2.Cells.Item(20,3).Value = "main"
$ws2.Range("A2:C2").Copy()
$ws2.Range("A20:C20").PasteSpecial(-4122)
$ws2.Cells.Item(20,4).Value = 19
$ws2.Cells.Item(20,5).Value = "Kierstein"
$ws2.Cells.Item(20,6).Value = "Failed"
$ws2.Cells.Item(20,7).Value = "Der Benutzer soll nach einem Fehler beim Suchen wieder Input geben können, aber ich kann keinen Fehler erzwingen"

# --- sheet1 row 21 -----------------------------------------------------------
$ws1.Cells.Item(21,1).Value = 20
$ws1.Cells.Item(21,2).Value = "7, Search"
$ws1.Cells.Item(21,3).Value = "Alternative Flow 4"
$ws1.Cells.Item(21,4).Value = "Benutzerinput ist leer und mit OK bestätigen"
$ws1.Cells.Item(21,5).Value = "Test database(Dokumente vorhanden)"
$ws1.Cells.Item(21,6).Value = "Alle vorhandenen Dokumente werden angezeigt"
$ws1.Range("E2:F2").Copy()
$ws1.Range("E21:F21").PasteSpecial(-4122)
$ws1.Rows.Item(21).RowHeight = 31.5

# --- sheet1 row 22 (G filled in before the rest, matching author's order) --
$ws1.Cells.Item(22,7).Value = "Benutzerinput ist identisch mit einem Dokumentenamen"
$ws1.Cells.Item(22,1).Value = 21
$ws1.Cells.Item(22,2).Value = "7, Search"
$ws1.Cells.Item(22,3).Value = "Alternative Flow 5"
$ws1.Cells.Item(22,4).Value = "Benutzerinput cdew.pdf und mit OK bestätigen"
$ws1.Cells.Item(22,5).Value = "Test database(cdew.pdf vorhanden)"
$ws1.Cells.Item(22,6).Value = "Meldung enthält nur das Dokument cdew.pdf"
$ws1.Range("E2:F2").Copy()
$ws1.Range("E22:F22").PasteSpecial(-4122)
$ws1.Rows.Item(22).RowHeight = 31.5

# --- sheet2 rows 21 / 22 -----------------------------------------------------
$ws2.Cells.Item(21,1).Value = 45132
$ws2.Cells.Item(21,2).Value = 0.69513888888888886
$ws2.Cells.Item(21,3).Value = "main"
$ws2.Range("A2:C2").Copy()
$ws2.Range("A21:C21").PasteSpecial(-4122)
$ws2.Cells.Item(21,4).Value = 20
$ws2.Cells.Item(21,5).Value = "Kierstein"
$ws2.Cells.Item(21,6).Value = "Success"

$ws2.Cells.Item(22,1).Value = 45132
$ws2.Cells.Item(22,2).Value = 0.6972222222222223
$ws2.Cells.Item(22,3).Value = "main"
$ws2.Range("A2:C2").Copy()
$ws2.Range("A22:C22").PasteSpecial(-4122)
$ws2.Cells.Item(22,4).Value = 21
$ws2.Cells.Item(22,5).Value = "Kierstein"
$ws2.Cells.Item(22,6).Value = "Success"

# ---------------------------------------------------------------------------
# Extend the Failed/Success conditional formatting to the new G19:G20 cells
# ---------------------------------------------------------------------------
$cfRange = $ws2.Range("G19:G20")
$cfFailed = $cfRange.FormatConditions.Add(9, 0, "Failed")
$cfFailed.Text = "Failed"
$cfFailed.Formula1 = '=NOT(ISERROR(SEARCH("Failed",F1)))'
$cfFailed.Interior.Color = 13551615
$cfFailed.Font.Color = 393372
$cfSuccess = $cfRange.FormatConditions.Add(9, 0, "Success")
$cfSuccess.Text = "Success"
$cfSuccess.Formula1 = '=NOT(ISERROR(SEARCH("Success",F1)))'
$cfSuccess.Interior.Color = 13561798
$cfSuccess.Font.Color = 24832

# ---------------------------------------------------------------------------
# Selection / active sheet: the workbook was left with sheet1 selected at
# B23 and sheet2 (now the active tab) selected at A23.
# ---------------------------------------------------------------------------
$ws1.Range("B23").Select()
$ws2.Activate()
$ws2.Range("A23").Select()
